$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New matchup rows to append (spring 23 week 7 inputs)
$data = @(
    @(4,0,5,3),
    @(2,0,6,3),
    @(3,2,5,0),
    @(7,3,5,0),
    @(5,3,2,0),
    @(6,3,5,0),
    @(3,1,3,2),
    @(4,0,5,2),
    @(3,0,3,3),
    @(7,2,6,0),
    @(3,3,4,0),
    @(3,1,3,2),
    @(5,2,5,0),
    @(5,2,5,0),
    @(5,2,4,0),
    @(5,3,3,0),
    @(5,2,5,1),
    @(6,2,6,0),
    @(4,0,4,2),
    @(3,2,4,1),
    @(4,3,3,0),
    @(2,0,3,3),
    @(5,2,5,0),
    @(5,3,5,0),
    @(6,0,5,2)
)

$startRow = 1687

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$lastRow = $startRow + $data.Length - 1
$nextRow = $lastRow + 1
$ws.Range("A$nextRow").Select()
$excel.ActiveWindow.ScrollRow = 1683
